$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2015")
$ws2016 = $wb.Worksheets.Item("2016")

# Column L header: copy K1's format (border/general), then set the text.
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)
$ws.Range("L1").Value = "Violent Crime Rate per 1,000 Residents"

# Column L sub-header: the "vcr1" field id uses the same format already
# present on the 2016 sheet's J2 cell (plain bordered style), not K2's.
$ws2016.Range("J2").Copy()
$ws.Range("L2").PasteSpecial(-4122)
$ws.Range("L2").Value = "vcr1"

# Data rows: copy each row's K-column format into L, then set the value.
$values = @(3.5, 1.9, 0.9, 4.18, 1.49, 4.6500000000000004, 1.46, 5.9, 1.63, 1.6, 4.71)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 11).Copy()
    $ws.Cells.Item($row, 12).PasteSpecial(-4122)
    $ws.Cells.Item($row, 12).Value = $values[$i]
}

$ws.Range("K16").Select()
